$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (drop bestFit, set explicit custom width)
$ws.Columns.Item(3).ColumnWidth = 11.666666666666666
$ws.Columns.Item(4).ColumnWidth = 44

# Row heights
$ws.Rows.Item(2).RowHeight = 57.6
$ws.Rows.Item(3).RowHeight = 43.2
$ws.Rows.Item(4).RowHeight = 43.2

# Row 2
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = "Sample Scene"
$ws.Range("D2").Value = "Run the game in a 16:9 screen aspect ratio (on maximum play as well). `nObserve any changes to the user interface across different screen sizes on a 16:9 aspect ratio."
$ws.Range("E2").Value = "The UI should not warp or alter on a 16:9 screen aspect ratio"

# Row 3
$ws.Range("D3").Value = "Run the game in a 16:9 screen aspect ratio.`nCollide the player with different objects and observe the health bar's changes"
$ws.Range("E3").Value = "The UI for the health bar at the top should indicate that the player is losing health (i.e. it starts out green at maximum health and shows more red as it loses health)"

# Row 4
$ws.Range("B4").Value = "6(b)"
$ws.Range("D4").Value = "Run the game in a 16:9 screen aspect ratio. `nObserve the timer as it runs in real time from seconds to minutes of starting the game."
$ws.Range("B3").Value = "6(a)"
$ws.Range("E4").Value = "The UI for the timer should be displayed in a real time '00:00:00' in minutes, seconds and milliseconds- in that order. The timer should keep "

# Row 3 and 4 also need "Sample Scene" in column C
$ws.Range("C3").Value = "Sample Scene"
$ws.Range("C4").Value = "Sample Scene"

# Alignment / wrap formatting
$ws.Range("D2:D4").HorizontalAlignment = -4131
$ws.Range("D2:D4").VerticalAlignment = -4160
$ws.Range("D2:D4").WrapText = $true

$ws.Range("E2").VerticalAlignment = -4160
$ws.Range("E3").WrapText = $true

# Selection
[void]$ws.Range("E5").Select()

Write-Host "done"
